$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last_update timestamp for the bevnat_info row (row 5, column E)
$ws.Range("E5").Value = 1706219962

# Update the active cell / selection on the sheet to E9
$null = $ws.Range("E9").Select()
